$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B16 value changes from 2 to 20
$ws.Range("B16").Value = 20

# New summary rows (row 20 stays empty)
$ws.Range("B21").Formula = "=SUM(B5:B13)"
$ws.Range("C21").Value = "Bartek"

$ws.Range("B22").Formula = "=SUM(B2:B4)"
$ws.Range("C22").Value = "Igor"

$ws.Range("B23").Formula = "=SUM(B14:B19)"
$ws.Range("C23").Value = "Justyna"

$ws.Range("B24").Formula = "=B21+B22+B23"
$ws.Range("C24").Value = "zespół"

# Bold the grand-total row
$ws.Range("B24:C24").Font.Bold = $true

# Match the selection left after the edit
$ws.Range("B24:C24").Select()
